$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A:F to 16 characters (ColumnWidth uses a slightly smaller
# "raw" unit than the stored character width, so 15.15 round-trips to 16).
$ws.Range("A1:F1").EntireColumn.ColumnWidth = 15.15

# Replace the abbreviated disease-name headers with their full names
# (typed in this order so the shared-strings table is appended to the
# same way the original author's edit produced it).
$ws.Range("B1").Value = "Amyotrophic lateral sclerosis"
$ws.Range("E1").Value = "Parkinson's disease"
$ws.Range("A1").Value = "Alzheimer's disease"
$ws.Range("C1").Value = "Dementia with`nLewy Bodies"
$ws.Range("D1").Value = "Frontotemporal`ndementia"

# The two multi-line headers wrap.
$ws.Range("C1:D1").WrapText = $true

# Row grows to fit the wrapped, multi-line headers.
$ws.Rows.Item(1).RowHeight = 43.2

# Move the active selection.
$ws.Range("I7").Select()
